$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Values -------------------------------------------------------
$ws.Range("A2").Value = "ซีพียู"
$ws.Range("B2").Value = 7000
$ws.Range("C2").Value = 0

$ws.Range("A3").Value = "เมนบอร์ด"
$ws.Range("B3").Value = 6000
$ws.Range("C3").Value = 0

$ws.Range("A4").Value = "แรม"
$ws.Range("B4").Value = 3500
$ws.Range("C4").Value = 0

$ws.Range("A5").Value = "โคมไฟตั้งพื้น"
$ws.Range("B5").Value = 2500
$ws.Range("C5").Value = 1
